$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 (an orphan row holding only B13/C13 = "5817535 - Lucas Barboza
# Sarno da Silva", with no label in column A) is removed. This shifts every
# row below it up by one (old row 14 -> new row 13, ... old row 25 -> new row 24).
$ws.Rows.Item(13).Delete()

# After the shift, a handful of B/C cells were overwritten with new text
# (rather than keeping the text that shifted into them). Apply those fixups.
$ws.Range("B10").Value = "5817535 - Lucas Barboza Sarno da Silva"
$ws.Range("C10").Value = "5817535 - Lucas Barboza Sarno da Silva"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

$ws.Range("B18").Value = "5817535 - Lucas Barboza Sarno da Silva"
$ws.Range("C18").Value = "5817535 - Lucas Barboza Sarno da Silva"

$nfMethodText = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("B19").Value = $nfMethodText
$ws.Range("C19").Value = $nfMethodText

$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

$rcText = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("B21").Value = $rcText
$ws.Range("C21").Value = $rcText
